# Auto-generated edit script: updates crypto price/volume table
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '55.814.79'
$ws.Range('E2').Value = '  +2.82%  '
$ws.Range('D3').Value = '2.506.57'
$ws.Range('E3').Value = '  +9.23%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '481.81'
$ws.Range('E5').Value = '  +6.65%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.16'
$ws.Range('E6').Value = '  +7.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.60%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.510'
$ws.Range('E8').Value = '  +7.14%  '
$ws.Range('D9').Value = '2.506.40'
$ws.Range('E9').Value = '  +10.09%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0986'
$ws.Range('E10').Value = '  +6.14%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.45'
$ws.Range('E11').Value = '  +0.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.327'
$ws.Range('E12').Value = '  +4.51%  '
$ws.Range('E13').Value = '  -0.15%  '
$ws.Range('D14').Value = '2.940.07'
$ws.Range('E14').Value = '  +9.80%  '
$ws.Range('D15').Value = '55.843.37'
$ws.Range('E15').Value = '  +2.93%  '
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.0000137'
$ws.Range('E16').Value = '  +13.15%  '
$ws.Range('B17').Value = 'Avalanche'
$ws.Range('C17').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '20.44'
$ws.Range('E17').Value = '  +7.79%  '
$ws.Range('D18').Value = '2.500.75'
$ws.Range('E18').Value = '  +10.24%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.34'
$ws.Range('E19').Value = '  +5.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '320.81'
$ws.Range('E20').Value = '  +5.10%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.95'
$ws.Range('E21').Value = '  +4.08%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.00'
$ws.Range('E22').Value = '  +0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.67'
$ws.Range('E23').Value = '  +5.40%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '57.84'
$ws.Range('E24').Value = '  +3.35%  '
$ws.Range('E25').Value = '  +3.60%  '
$ws.Range('E26').Value = '  +2.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.403'
$ws.Range('E27').Value = '  +6.90%  '
$ws.Range('D28').Value = '2.614.11'
$ws.Range('E28').Value = '  +10.59%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.40'
$ws.Range('E29').Value = '  +7.14%  '
$ws.Range('D30').Value = '0.0₃0766'
$ws.Range('E30').Value = '  +7.28%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('E31').Value = '  +0.69%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '149.85'
$ws.Range('E32').Value = '  +3.25%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '18.11'
$ws.Range('E33').Value = '  +6.01%  '
$ws.Range('E34').Value = '  +8.83%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.18'
$ws.Range('E35').Value = '  +8.30%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.69'
$ws.Range('E36').Value = '  +0.67%  '
$ws.Range('E37').Value = '  +8.46%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.844'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '34.51'
$ws.Range('E39').Value = '  +4.30%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.614'
$ws.Range('E40').Value = '  +19.77%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.999'
$ws.Range('E41').Value = '  +0.59%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.33'
$ws.Range('E42').Value = '  +5.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0545'
$ws.Range('E43').Value = '  +8.15%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '3.37'
$ws.Range('E44').Value = '  +5.65%  '
$ws.Range('B45').Value = 'WhiteBITCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.15'
$ws.Range('E45').Value = '  -1.60%  '
$ws.Range('B46').Value = 'Maker'
$ws.Range('C46').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D46').Value = '1.969.89'
$ws.Range('E46').Value = '  +1.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0901'
$ws.Range('E47').Value = '  +10.14%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '251.43'
$ws.Range('E48').Value = '  +30.93%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0222'
$ws.Range('E49').Value = '  +6.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '17.50'
$ws.Range('E50').Value = '  +5.91%  '
$ws.Range('E51').Value = '  +6.34%  '
